$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 62.666668
$ws.Range("I6").Value = 62.666668
$ws.Range("K6").Value = 188.000004
$ws.Range("M6").Value = -76.00000399999999
# Row 33
$ws.Range("H33").Value = 7815323.5
$ws.Range("I33").Value = 8928945
$ws.Range("J33").Value = 19975
$ws.Range("K33").Value = 8928945
$ws.Range("L33").Value = 19975
$ws.Range("M33").Value = -8928716
$ws.Range("N33").Value = -20433
# Row 62
$ws.Range("H62").Value = 5742.577
$ws.Range("I62").Value = 3016.9473
$ws.Range("K62").Value = 3016.9473
$ws.Range("M62").Value = -2392.9473
# Row 65
$ws.Range("H65").Value = 5742.577
$ws.Range("I65").Value = 3016.9473
$ws.Range("K65").Value = 15084.7365
$ws.Range("M65").Value = -11964.7365
# Row 127
$ws.Range("H127").Value = 701.8461
$ws.Range("I127").Value = 465.81818
$ws.Range("K127").Value = 1397.45454
$ws.Range("M127").Value = 3562.54546
# Row 132
$ws.Range("H132").Value = 6715.551
$ws.Range("I132").Value = 5362.6587
$ws.Range("J132").Value = 13649.125
$ws.Range("K132").Value = 16087.9761
$ws.Range("L132").Value = 40947.375
$ws.Range("M132").Value = -13557.9761
$ws.Range("N132").Value = -46007.375
# Row 138
$ws.Range("H138").Value = 2307.14
$ws.Range("I138").Value = 1181.8
$ws.Range("J138").Value = 2588.475
$ws.Range("K138").Value = 3545.4
$ws.Range("L138").Value = 7765.424999999999
$ws.Range("M138").Value = 1594.6
$ws.Range("N138").Value = -18045.425

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6350.34
$ws.Range("I32").Value = 1154.7028
$ws.Range("J32").Value = 21137.924
$ws.Range("K32").Value = 1154.7028
$ws.Range("L32").Value = 21137.924
$ws.Range("M32").Value = -867.7028
$ws.Range("N32").Value = -21711.924
# Row 74
$ws.Range("H74").Value = 10699.305
$ws.Range("I74").Value = 1852.5938
$ws.Range("J74").Value = 30920.357
$ws.Range("K74").Value = 1852.5938
$ws.Range("L74").Value = 30920.357
$ws.Range("M74").Value = -978.5938000000001
$ws.Range("N74").Value = -32668.357
# Row 77
$ws.Range("H77").Value = 10699.305
$ws.Range("I77").Value = 1852.5938
$ws.Range("J77").Value = 30920.357
$ws.Range("K77").Value = 9262.969000000001
$ws.Range("L77").Value = 154601.785
$ws.Range("M77").Value = -4894.969000000001
$ws.Range("N77").Value = -163337.785

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 19771.564
$ws.Range("I20").Value = 7133.9473
$ws.Range("J20").Value = 28664.703
$ws.Range("K20").Value = 7133.9473
$ws.Range("L20").Value = 28664.703
$ws.Range("M20").Value = -6886.9473
$ws.Range("N20").Value = -29158.703
# Row 107
$ws.Range("H107").Value = 1513
$ws.Range("I107").Value = 1018.2
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 1018.2
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = 901.8
$ws.Range("N107").Value = -6590

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5932.8
$ws.Range("I16").Value = 918.4286
$ws.Range("J16").Value = 17633
$ws.Range("K16").Value = 918.4286
$ws.Range("L16").Value = 17633
$ws.Range("M16").Value = -631.4286
$ws.Range("N16").Value = -18207
# Row 31
$ws.Range("H31").Value = 10308.462
$ws.Range("I31").Value = 4997.759
$ws.Range("J31").Value = 17004.564
$ws.Range("K31").Value = 4997.759
$ws.Range("L31").Value = 17004.564
$ws.Range("M31").Value = -4702.759
$ws.Range("N31").Value = -17594.564
# Row 32
$ws.Range("H32").Value = 2288.5
$ws.Range("I32").Value = 3339.6667
$ws.Range("J32").Value = 1237.3334
$ws.Range("K32").Value = 3339.6667
$ws.Range("L32").Value = 1237.3334
$ws.Range("M32").Value = -3023.6667
$ws.Range("N32").Value = -1869.3334
# Row 34
$ws.Range("H34").Value = 10308.462
$ws.Range("I34").Value = 4997.759
$ws.Range("J34").Value = 17004.564
$ws.Range("K34").Value = 4997.759
$ws.Range("L34").Value = 17004.564
$ws.Range("M34").Value = -4795.759
$ws.Range("N34").Value = -17408.564
# Row 58
$ws.Range("H58").Value = 9558.883
$ws.Range("I58").Value = 3768.5518
$ws.Range("J58").Value = 17191.592
$ws.Range("K58").Value = 3768.5518
$ws.Range("L58").Value = 17191.592
$ws.Range("M58").Value = -3565.5518
$ws.Range("N58").Value = -17597.592
# Row 86
$ws.Range("H86").Value = 8380.261
$ws.Range("I86").Value = 9913.416999999999
$ws.Range("J86").Value = 6707.727
$ws.Range("K86").Value = 9913.416999999999
$ws.Range("L86").Value = 6707.727
$ws.Range("M86").Value = -8790.416999999999
$ws.Range("N86").Value = -8953.726999999999
# Row 87
$ws.Range("H87").Value = 37400
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 37400
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 37400
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -39772
# Row 89
$ws.Range("H89").Value = 8380.261
$ws.Range("I89").Value = 9913.416999999999
$ws.Range("J89").Value = 6707.727
$ws.Range("K89").Value = 49567.085
$ws.Range("L89").Value = 33538.635
$ws.Range("M89").Value = -43951.085
$ws.Range("N89").Value = -44770.635
# Row 90
$ws.Range("H90").Value = 37400
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 37400
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 112200
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -124056
# Row 113
$ws.Range("H113").Value = 5932.8
$ws.Range("I113").Value = 918.4286
$ws.Range("J113").Value = 17633
$ws.Range("K113").Value = 918.4286
$ws.Range("L113").Value = 17633
$ws.Range("M113").Value = 1251.5714
$ws.Range("N113").Value = -21973
# Row 132
$ws.Range("H132").Value = 3424.8867
$ws.Range("I132").Value = 1017.2955
$ws.Range("J132").Value = 15195.333
$ws.Range("K132").Value = 3051.8865
$ws.Range("L132").Value = 45585.999
$ws.Range("M132").Value = -521.8864999999996
$ws.Range("N132").Value = -50645.999
# Row 134
$ws.Range("H134").Value = 23261058
$ws.Range("I134").Value = 1481.8334
$ws.Range("J134").Value = 52641576
$ws.Range("K134").Value = 4445.5002
$ws.Range("L134").Value = 157924728
$ws.Range("M134").Value = -1910.5002
$ws.Range("N134").Value = -157929798
# Row 136
$ws.Range("H136").Value = 9558.883
$ws.Range("I136").Value = 3768.5518
$ws.Range("J136").Value = 17191.592
$ws.Range("K136").Value = 11305.6554
$ws.Range("L136").Value = 51574.776
$ws.Range("M136").Value = -8755.6554
$ws.Range("N136").Value = -56674.776

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 13457247
$ws.Range("I122").Value = 23359504
$ws.Range("J122").Value = 3554989.5
$ws.Range("K122").Value = 210235536
$ws.Range("L122").Value = 31994905.5
$ws.Range("M122").Value = -210233086
$ws.Range("N122").Value = -31999805.5
# Row 139
$ws.Range("H139").Value = 12633.579
$ws.Range("I139").Value = 69069
$ws.Range("J139").Value = 5994.1177
$ws.Range("K139").Value = 207207
$ws.Range("L139").Value = 17982.3531
$ws.Range("M139").Value = -202067
$ws.Range("N139").Value = -28262.3531

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 694920.9
$ws.Range("I132").Value = 1314.5333
$ws.Range("J132").Value = 3095866
$ws.Range("K132").Value = 3943.5999
$ws.Range("L132").Value = 9287598
$ws.Range("M132").Value = -1413.5999
$ws.Range("N132").Value = -9292658
# Row 139
$ws.Range("H139").Value = 69715
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

